$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.478.59'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.71%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.301.76'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.10%  '

$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.13%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.51'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.13%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.630'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.45%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.607'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.16%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.71'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.76%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0910'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.53%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.36'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.75%  '

$ws.Range('E13').Value = '  +0.29%  '

$ws.Range('E14').Value = '  -1.54%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.25'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.28%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.651.24'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.08%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.303.26'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.06%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.440.26'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.52%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.47'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.31%  '

$ws.Range('E20').Value = '  +0.83%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.79%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.53'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.18%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '276.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +6.78%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.29'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +20.75%  '

$ws.Range('E25').Value = '  -1.24%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.84'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.54%  '

$ws.Range('E28').Value = '  +3.04%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '22.76'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.41%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.82'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.56%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '165.49'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.08%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0874'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.73%  '

$ws.Range('E33').Value = '  +0.32%  '

$ws.Range('E34').Value = '  +4.06%  '

$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.61'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -10.95%  '

$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.118'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.84%  '

$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0367'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.55%  '

$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.57'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.03%  '

$ws.Range('E40').Value = '  -0.15%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.51'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.84%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '69.66'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.98%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '94.73'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.33%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.227'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.72%  '

$ws.Range('E45').Value = '  -0.12%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '81.17'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +7.93%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.07'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.49%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '113.04'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.49%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.97'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.00%  '

$ws.Range('E50').Value = '  -2.33%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.590.27'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.04%  '
